# Automatische test-sync: 2025-08-28 17:37:50
# Appends a new "Retour status" log row to the Logs sheet, extends the
# conditional-formatting ranges to cover it, and updates the Dashboard
# summary count for "Retour / Terugbetaling".

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: add row 3 ----
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Retour status"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("D3").Value = "Retour / Terugbetaling"
$logs.Range("F3").Value = "2025-08-28 17:37:41"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Nee"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# ---- Extend conditional formatting ranges from row 2 to rows 2:3 ----
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range($col + "2")
    $newRange = $logs.Range($col + "2:" + $col + "3")
    $fc = $oldRange.FormatConditions
    for ($i = 1; $i -le $fc.Count(); $i++) {
        $fc.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---- Dashboard sheet: bump the Retour / Terugbetaling count ----
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 2
